$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "P14" = 8.86666666666667
    "P16" = 7.26666666666667
    "P17" = 6.66666666666667
    "P20" = 5.46666666666667
    "P24" = 4.23333333333333
    "P26" = 3.96666666666667
    "P28" = 4.36666666666667
    "P30" = 4.53333333333333
    "P31" = 4.63333333333333
    "P33" = 5.26666666666667
    "P34" = 4.83333333333333
    "P37" = 4.76666666666667
    "P38" = 5.13333333333333
    "P39" = 4.66666666666667
    "P42" = 4.43333333333333
    "P43" = 4.86666666666667
    "P45" = 4.53333333333333
    "P47" = 4.73333333333333
    "P49" = 4.63333333333333
    "P53" = 5.23333333333333
    "P54" = 5.36666666666667
    "P55" = 5.93333333333333
    "P56" = 7.36666666666667
    "P57" = 8.53333333333333
    "P58" = 10.8666666666667
    "P61" = 13.8333333333333
    "P63" = 14.3666666666667
    "P64" = 14.6333333333333
    "P65" = 15.6333333333333
    "P66" = 15.1666666666667
    "P67" = 15.0333333333333
    "P69" = 15.8333333333333
    "P70" = 15.9333333333333
    "P71" = 15.5333333333333
    "P74" = 14.5666666666667
    "P75" = 14.3666666666667
    "P76" = 13.3333333333333
    "P77" = 12.8333333333333
    "P78" = 12.8333333333333
    "P80" = 11.6666666666667
    "P81" = 10.8666666666667
    "P82" = 10.6333333333333
    "P84" = 9.53333333333333
    "P85" = 9.43333333333333
    "P86" = 8.93333333333333
    "P88" = 8.23333333333333
    "P89" = 7.53333333333333
    "P90" = 7.26666666666667
    "P91" = 6.66666666666667
    "P92" = 6.66666666666667
    "P93" = 6.43333333333333
    "P94" = 5.96666666666667
    "P95" = 5.86666666666667
    "P99" = 5.16666666666667
    "P102" = 4.86666666666667
    "P103" = 5.13333333333333
    "P105" = 6.26666666666667
    "P106" = 7.33333333333333
    "P107" = 6.96666666666667
    "P108" = 5.46666666666667
    "M109" = 19.2
    "P109" = 5.23333333333333
    "P110" = 4.86666666666667
    "P112" = 4.26666666666667
    "P113" = 4.43333333333333
    "P115" = 4.16666666666667
    "M116" = 13.3
    "P116" = 4.43333333333333
    "M117" = 12.9
    "P117" = 4.43333333333333
    "I118" = 42705
    "M118" = 14.1
    "O118" = 2452566.66666667
    "P118" = 4.26666666666667
    "I119" = 43050
    "O119" = 2492200
    "P119" = 4.36666666666667
    "S119" = 730.42
    "T119" = 54221130
    "F120" = 156.7
    "I120" = 43825
    "M120" = 14.3
    "O120" = 2506966.66666667
    "P120" = 4.26666666666667
    "C121" = 3894691
    "D121" = 116.733333333333
    "E121" = 23263588708.15
    "L121" = 119.566666666667
    "P121" = 4.3
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
